# New crime data collected - weekly CompStat refresh for 66th Precinct
# Updates: report volume/date header text, plus the Week/28-Day/YTD/2-Year
# crime-complaint figures (and their derived % changes) in rows 15-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings) - volume number and report week
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/18/2023  Through  9/24/2023"

# ---------------------------------------------------------------------
# Row 15 - Rape : Week-to-date 2023 count drops to 0 (becomes the "0"
# text placeholder, matching the style already used by sibling cells).
# ---------------------------------------------------------------------
$ws.Range("D15").Copy($ws.Range("C15"))

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 63
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = -5.970149253731
$ws.Range("L16").Value = 18.867924528301
$ws.Range("M16").Value = -51.162790697674
$ws.Range("N16").Value = -87.766990291262

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 154
$ws.Range("J17").Value = 166
$ws.Range("K17").Value = -7.228915662650
$ws.Range("L17").Value = 31.623931623931
$ws.Range("M17").Value = 31.623931623931
$ws.Range("N17").Value = -36.625514403292

# ---------------------------------------------------------------------
# Row 18 - Burglary : Week-to-date 2023 count goes from 0 to an actual
# number, so C18 flips from the "0" text placeholder to a real number.
# ---------------------------------------------------------------------
$ws.Range("I18").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 85
$ws.Range("J18").Value = 139
$ws.Range("K18").Value = -38.848920863309
$ws.Range("L18").Value = -22.018348623853
$ws.Range("M18").Value = -69.642857142857
$ws.Range("N18").Value = -93.227091633466

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("E19").Value = -30.769230769230
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 15.909090909090
$ws.Range("I19").Value = 420
$ws.Range("J19").Value = 454
$ws.Range("K19").Value = -7.488986784140
$ws.Range("L19").Value = 33.757961783439
$ws.Range("M19").Value = 46.341463414634
$ws.Range("N19").Value = -11.949685534591

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 7.692307692307
$ws.Range("I20").Value = 115
$ws.Range("J20").Value = 84
$ws.Range("K20").Value = 36.904761904761
$ws.Range("L20").Value = 94.915254237288
$ws.Range("M20").Value = 8.490566037735
$ws.Range("N20").Value = -91.301059001512

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -4.761904761904
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = -5.769230769230
$ws.Range("I21").Value = 853
$ws.Range("J21").Value = 925
$ws.Range("K21").Value = -7.783783783783
$ws.Range("L21").Value = 29.242424242424
$ws.Range("M21").Value = -7.883369330453
$ws.Range("N21").Value = -77.751695357329

# ---------------------------------------------------------------------
# Row 22 - Transit : Week-to-date 2023 counts (C22/F22) go from 0 to
# real numbers, so those two flip from "0" text placeholders to numbers.
# ---------------------------------------------------------------------
$ws.Range("J22").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("J22").Copy($ws.Range("F22"))
$ws.Range("F22").Value = 1
$ws.Range("I22").Value = 8
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = -11.111111111111
$ws.Range("M22").Value = -50

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -40.740740740740
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -33.333333333333
$ws.Range("I24").Value = 797
$ws.Range("J24").Value = 838
$ws.Range("K24").Value = -4.892601431980
$ws.Range("L24").Value = 36.006825938566
$ws.Range("M24").Value = 18.777943368107

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 80
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -26.190476190476
$ws.Range("I25").Value = 256
$ws.Range("J25").Value = 255
$ws.Range("K25").Value = 0.392156862745
$ws.Range("L25").Value = 15.315315315315
$ws.Range("M25").Value = -19.242902208201

# ---------------------------------------------------------------------
# Row 26 - UCR Rape* : both 2023 and 2022 week-to-date counts drop to 0,
# so C26/D26 flip to the "0" text placeholder and E26's 0% chg becomes
# the "***.*" placeholder (division by zero).
# ---------------------------------------------------------------------
$ws.Range("C23").Copy($ws.Range("C26"))
$ws.Range("D23").Copy($ws.Range("D26"))
$ws.Range("E23").Copy($ws.Range("E26"))
$ws.Range("L26").Value = 33.333333333333

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes : 2022 week-to-date count drops to 0, so
# D27 flips to the "0" text placeholder and E27 becomes "***.*".
# ---------------------------------------------------------------------
$ws.Range("D28").Copy($ws.Range("D27"))
$ws.Range("E28").Copy($ws.Range("E27"))
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -66.666666666666
$ws.Range("L27").Value = -7.547169811320

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes : 2022 week-to-date count drops to 0, so D30
# flips to the "0" text placeholder and E30 becomes "***.*".
# ---------------------------------------------------------------------
$ws.Range("D29").Copy($ws.Range("D30"))
$ws.Range("E29").Copy($ws.Range("E30"))
$ws.Range("L30").Value = -33.333333333333
